# The workbook's single data row (row 2) had its "Pdf_URL" value (cell AL2,
# "http://arpeissig.at/wp-content/uploads/2016/02/D7_NHB_ARP_Final_2.pdf")
# removed/cleared, leaving that cell blank. The trailing "Report Html
# Address" / "Report Drt Address" / "Database link" cells (AM2:AO2) are
# untouched content-wise.
#
# The author also scrolled the sheet so column S is the left-most visible
# column and left the now-empty AL2 cell selected - reproduce that
# navigation as closely as the object model allows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll so column S is the first visible column (matches the saved view),
# then select and clear the Pdf_URL cell for the one data row.
$ws.Range("S1").Select()
$ws.Range("AL2").Select()
$ws.Range("AL2").ClearContents()
